# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.208.03"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.87"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.66"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5071"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3917"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09310"
$ws.Range("E9").Value = "  -3.40%  "

$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.87"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.388"
$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.92"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.903.79"
$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.314"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001123"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.46"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06607"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.267.02"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.43"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.593"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.126.96"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.96"
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.10"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.099"
$ws.Range("E31").Value = "  +2.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1075"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.631"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.639"
$ws.Range("E35").Value = "  +0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06665"
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.247"
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("E40").Value = "  +7.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6447"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.998"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6034"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.282"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.00"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.184"
$ws.Range("E51").Value = "  -1.02%  "
